# The deck's design was switched from the custom "Integral" theme to the
# built-in "Office Theme" palette (dark1/light1/dark2/light2/accent1-6/
# hyperlink/followed-hyperlink). Re-point the deck's theme color scheme at
# the standard Office Theme RGB values via the Design > Variants > Colors
# object model (PowerPoint.ThemeColorScheme).

# PowerPoint/VBA's RGB() packs bytes as (B*65536 + G*256 + R); rebuild that
# here since this host doesn't expose the VBA RGB() helper function.
function OfficeRGB($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme reaches the single color scheme shared by the
# whole deck (slide master + every layout), so slide 1 is as good an anchor
# as any other slide.
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Item() order follows the msoThemeColorXxx enumeration:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5-10 Accent1-6, 11 Hyperlink,
# 12 FollowedHyperlink.
$cs.Item(1).RGB  = OfficeRGB 0x00 0x00 0x00   # Dark 1             000000
$cs.Item(2).RGB  = OfficeRGB 0xFF 0xFF 0xFF   # Light 1            FFFFFF
$cs.Item(3).RGB  = OfficeRGB 0x44 0x54 0x6A   # Dark 2             44546A
$cs.Item(4).RGB  = OfficeRGB 0xE7 0xE6 0xE6   # Light 2            E7E6E6
$cs.Item(5).RGB  = OfficeRGB 0x5B 0x9B 0xD5   # Accent 1           5B9BD5
$cs.Item(6).RGB  = OfficeRGB 0xED 0x7D 0x31   # Accent 2           ED7D31
$cs.Item(7).RGB  = OfficeRGB 0xA5 0xA5 0xA5   # Accent 3           A5A5A5
$cs.Item(8).RGB  = OfficeRGB 0xFF 0xC0 0x00   # Accent 4           FFC000
$cs.Item(9).RGB  = OfficeRGB 0x44 0x72 0xC4   # Accent 5           4472C4
$cs.Item(10).RGB = OfficeRGB 0x70 0xAD 0x47   # Accent 6           70AD47
$cs.Item(11).RGB = OfficeRGB 0x05 0x63 0xC1   # Hyperlink          0563C1
$cs.Item(12).RGB = OfficeRGB 0x95 0x4F 0x72   # Followed Hyperlink 954F72
